# Chapter 8 (Conclusion) end-of-chapter wording fixes.
$d = $word.ActiveDocument

# 1. "additional contributions are made to" -> "additional contributions to"
$d.Content.Find.Execute(
    "this thesis has made additional contributions are made to the fields",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "this thesis has made additional contributions to the fields", 2) | Out-Null

# 2. "the research has contributed value" -> "the research has already contributed value"
$d.Content.Find.Execute(
    "the research has contributed value to real-world industrial projects",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "the research has already contributed value to real-world industrial projects", 2) | Out-Null

# 3. Bold "sets the scene for an progressive and activist agenda" and extend the
#    following clause with "and reconfigure society to one where those
#    human-centric needs are better met".
$bold = $d.Content.Duplicate
$bold.Find.Execute(
    "sets the scene for an progressive and activist agenda",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bold.Font.Bold = 1

$d.Content.Find.Execute(
    "to realise those needs. It constitutes",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "to realise those needs and reconfigure society to one where those human-centric needs are better met. It constitutes", 2) | Out-Null
